$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the B7/B8 row swap (Дария Златкова / Гергана Джонджорова) ---
# Row 7 currently holds "Гергана Джонджорова" / 1601681081 / 5
# Row 8 currently holds "Дария Златкова" / 1601681080 / 6
# Target: row 7 = Дария Златкова / 1601681080 / 6 ; row 8 = Гергана Джонджорова / 1601681081 / 5
$ws.Range("A7").Value = "Дария Златкова"
$ws.Range("B7").Value = 1601681080
$ws.Range("C7").Value = 6
$ws.Range("A8").Value = "Гергана Джонджорова"
$ws.Range("B8").Value = 1601681081
$ws.Range("C8").Value = 5

# --- Expand the table to include the two new "control" grade columns
#     and the running/current grade column ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E19"))

# Rename existing "Оценка" header to "Оценка контролно 1" and add the two
# new headers; writing into the table header row renames the ListColumns.
$ws.Range("C1").Value = "Оценка контролно 1"
$ws.Range("D1").Value = "Оценка контролно 2"
$ws.Range("E1").Value = "Текуща оценка"

# --- Fill in "Оценка контролно 2" (D) and "Текуща оценка" (E) values ---
$data = @(
  @{ row = 2;  d = "N/A"; e = 4 },
  @{ row = 3;  d = 6;     e = 6 },
  @{ row = 4;  d = 6;     e = 5 },
  @{ row = 5;  d = 6;     e = 6 },
  @{ row = 6;  d = 6;     e = 5 },
  @{ row = 7;  d = 6;     e = 6 },
  @{ row = 8;  d = 6;     e = 6 },
  @{ row = 9;  d = 6;     e = 5 },
  @{ row = 10; d = 6;     e = 5 },
  @{ row = 11; d = 6;     e = 5 },
  @{ row = 12; d = 6;     e = 5 },
  @{ row = 13; d = 6;     e = 6 },
  @{ row = 14; d = 6;     e = 6 },
  @{ row = 15; d = 6;     e = 6 },
  @{ row = 16; d = 6;     e = 6 },
  @{ row = 17; d = 6;     e = 6 },
  @{ row = 18; d = 6;     e = 6 }
)

foreach ($entry in $data) {
  $r = $entry.row
  $ws.Range("D$r").Value = $entry.d
  $ws.Range("E$r").Value = $entry.e
}

# D2 ("N/A") and D12 should be right aligned, matching the style added for
# the "N/A" placeholder cell.
$ws.Range("D2").HorizontalAlignment = -4152
$ws.Range("D12").HorizontalAlignment = -4152

# Widen the new/changed columns (values chosen so the stored OOXML width,
# which this engine rounds to the nearest 1/6 character, lands as close as
# possible to the author's original widths).
$ws.Columns.Item(3).ColumnWidth = 25.799479166666668
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668
$ws.Columns.Item(5).ColumnWidth = 16.709635416666668

# Move the active selection back to B1 (matches the saved view state).
$ws.Range("B1").Select()
